$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.524.88'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '2.470.00'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.45'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.82'
$ws.Range("E6").Value = '  -3.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.14'
$ws.Range("E10").Value = '  -4.60%  '
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").Value = '2.848.67'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.85'
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.02'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").Value = '2.431.65'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.768'
$ws.Range("E17").Value = '  -3.49%  '
$ws.Range("D18").Value = '41.496.33'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.49'
$ws.Range("E19").Value = '  +2.23%  '
$ws.Range("D20").Value = '0.0₃0946'
$ws.Range("E20").Value = '  +1.90%  '
$ws.Range("E21").Value = '  +3.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.10'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.02'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.72'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.89'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.55'
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.30'
$ws.Range("E30").Value = '  -3.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.10'
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.44'
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.20'
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.88'
$ws.Range("E36").Value = '  -6.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.104'
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.114'
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.79'
$ws.Range("E39").Value = '  -4.98%  '
$ws.Range("B40").Value = 'ApeXProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  -11.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.02'
$ws.Range("E41").Value = '  -5.27%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '1.940.46'
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0283'
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.39'
$ws.Range("E45").Value = '  -6.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.04'
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("D48").Value = '2.705.22'
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.03'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.05'
$ws.Range("E50").Value = '  -4.43%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.26'
$ws.Range("E51").Value = '  +2.33%  '
